# The "About" sheet's C1 cell held a one-off date stamp (value 44307,
# formatted with the date-only number format style). This edit removes
# that cell entirely (value + formatting), restoring the sheet to just
# the title in A1 for row 1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")
$ws.Range("C1").Clear()
